# Apply evaluation data edits for constellations up to number 38
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("B21").Value = 0.1
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 10
$ws.Range("I21").Value = 8
$ws.Range("J21").Value = "GreLum"

# Row 22
$ws.Range("B22").Value = 0.1
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = "GreLum"

# Row 23
$ws.Range("B23").Value = 0.1
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 10
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = "GreLum, LowDiff"

# Row 24
$ws.Range("B24").Value = 0.01
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 10
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = "GreLum, LowDiff"

# Row 25
$ws.Range("B25").Value = 0.05
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 10
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = "GreLum, LowDiff"

# Row 26
$ws.Range("B26").Value = 0.05
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 10
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = "GreLum, LowDiff"

# Row 27
$ws.Range("B27").Value = 0.5
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 20
$ws.Range("F27").Value = 6
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = "GreLum, LowDiff"

# Row 28
$ws.Range("B28").Value = 0.05
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 10
$ws.Range("I28").Value = 0

# Row 29
$ws.Range("B29").Value = 0.05
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 10
$ws.Range("I29").Value = 0

# Row 30
$ws.Range("B30").Value = 0.05
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 10
$ws.Range("I30").Value = 0

# Row 31
$ws.Range("B31").Value = 0.05
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 10
$ws.Range("I31").Value = 2
$ws.Range("J31").Value = "GreLum"

# Row 32
$ws.Range("B32").Value = 0.05
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 10
$ws.Range("I32").Value = 2
$ws.Range("J32").Value = "GreLum"

# Row 33
$ws.Range("B33").Value = 0.05
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 5
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 10
$ws.Range("I33").Value = 0

# Row 34
$ws.Range("B34").Value = 0.05
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 10
$ws.Range("I34").Value = 0

# Row 35
$ws.Range("B35").Value = 0.1
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 10
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 7
$ws.Range("I35").Value = 3
$ws.Range("J35").Value = "GreLum"

# Row 36
$ws.Range("J36").Value = "some cases GreLum, some LowDiff. Hard because of only 3 stars in constellation"

# Row 37
$ws.Range("B37").Value = 0.5
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 10
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 3
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = "GreLum"

# Row 39
$ws.Range("B39").Value = 0.03
$ws.Range("C39").Value = 3
$ws.Range("D39").Value = 4
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 10
$ws.Range("I39").Value = 4
$ws.Range("J39").Value = "GreLum"

# Row 40
$ws.Range("B40").Value = 0.05
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 5
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 10
$ws.Range("I40").Value = 0

# Recalculate formulas (H column = F/G) so cached values are correct
$excel.Calculate()

# Update the active sheet view/selection to match final state
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B41").Select()
